$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.635.34"
$ws.Range("E2").Value = "  -3.86%  "
$ws.Range("D3").Value = "2.972.51"
$ws.Range("E3").Value = "  -5.21%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.20"
$ws.Range("E5").Value = "  -4.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.22"
$ws.Range("E6").Value = "  -5.83%  "
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "2.982.06"
$ws.Range("E9").Value = "  -5.13%  "
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("E11").Value = "  -6.25%  "
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("D13").Value = "3.492.00"
$ws.Range("E13").Value = "  -5.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.124"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").Value = "61.708.89"
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.79"
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").Value = "2.982.21"
$ws.Range("E17").Value = "  -4.89%  "
$ws.Range("E18").Value = "  -4.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.18"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.04"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.77"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("E22").Value = "  -5.42%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.89"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").Value = "3.092.99"
$ws.Range("E26").Value = "  -5.76%  "
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0944"
$ws.Range("E28").Value = "  -6.17%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.33"
$ws.Range("E30").Value = "  -4.92%  "
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.50"
$ws.Range("E33").Value = "  -2.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "160.90"
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.67"
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.08"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("E39").Value = "  -5.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.93"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "2.414.19"
$ws.Range("E41").Value = "  -9.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.33"
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.29"
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.666"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.07"
$ws.Range("E48").Value = "  -6.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.95"
$ws.Range("E49").Value = "  -5.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "270.14"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0953"
$ws.Range("E51").Value = "  -2.05%  "
